# Add a second "Pipes" table that mirrors the existing "Sinais" table
# (columns D:E) one column gap to the right, in columns G:H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Merge the (currently blank) header cells G4:H4 first, while they are
#    still unstyled, so that the subsequent format-only paste below can
#    reuse the existing "Sinais" header style instead of Excel deriving a
#    brand-new split-border style pair for the merge.
$ws.Range("G4:H4").Merge()

# 2) Clone the look of the "Sinais" header (D4:E4) onto G4:H4 without
#    touching their (already-correct) merge state.
$ws.Range("D4:E4").Copy()
$ws.Range("G4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 3) Set the new header text via the shared-string table ("Pipes").
$ws.Range("G4").Value = "Pipes"

# 4) Mirror the 13 numbered rows of the "Sinais" table (D5:E17) into the
#    new "Pipes" table (G5:H17) in one shot, preserving values + styles.
$ws.Range("D5:E17").Copy($ws.Range("G5"))
$excel.CutCopyMode = 0

# 5) Restore the selection to G5, matching the saved workbook state.
$ws.Range("G5").Select()

# 6) Match the saved page setup (portrait, paper size 9 / A4).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
